# Updates the LR-pairs lrc2p table (Vtn -> Tnfrsf11b) with the new TPM-based
# NATMI run: adds the "Resolving-Mac" sending cluster and refreshes every
# derived statistic for the existing ECs/FAPs/MuSCs clusters (rows 2-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Vtn/Tnfrsf11b)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Vtn"
$ws.Cells.Item(2, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.767552
$ws.Cells.Item(2, 8).Value = 8.302655999999999
$ws.Cells.Item(2, 9).Value = 0.04706493447833917
$ws.Cells.Item(2, 10).Value = 0.04706493447833917
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.06861733333333334
$ws.Cells.Item(2, 14).Value = 0.205852
$ws.Cells.Item(2, 15).Value = 0.01654048691795588
$ws.Cells.Item(2, 16).Value = 0.01654048691795588
$ws.Cells.Item(2, 17).Value = 0.1899020381013333
$ws.Cells.Item(2, 18).Value = 1.709118342912
$ws.Cells.Item(2, 19).Value = 0.0007784769330334199
$ws.Cells.Item(2, 20).Value = 0.0007784769330334199

# Row 3: ECs -> FAPs (Vtn/Tnfrsf11b)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vtn"
$ws.Cells.Item(3, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.767552
$ws.Cells.Item(3, 8).Value = 8.302655999999999
$ws.Cells.Item(3, 9).Value = 0.04706493447833917
$ws.Cells.Item(3, 10).Value = 0.04706493447833917
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.776574666666666
$ws.Cells.Item(3, 14).Value = 11.329724
$ws.Cells.Item(3, 15).Value = 0.9103586635352137
$ws.Cells.Item(3, 16).Value = 0.9103586635352137
$ws.Cells.Item(3, 17).Value = 10.45186677188267
$ws.Cells.Item(3, 18).Value = 94.06680094694397
$ws.Cells.Item(3, 19).Value = 0.04284597085107324
$ws.Cells.Item(3, 20).Value = 0.04284597085107324

# Row 4: ECs -> MuSCs (Vtn/Tnfrsf11b)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Vtn"
$ws.Cells.Item(4, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.767552
$ws.Cells.Item(4, 8).Value = 8.302655999999999
$ws.Cells.Item(4, 9).Value = 0.04706493447833917
$ws.Cells.Item(4, 10).Value = 0.04706493447833917
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.303255
$ws.Cells.Item(4, 14).Value = 0.909765
$ws.Cells.Item(4, 15).Value = 0.07310084954683041
$ws.Cells.Item(4, 16).Value = 0.07310084954683042
$ws.Cells.Item(4, 17).Value = 0.8392739817599999
$ws.Cells.Item(4, 18).Value = 7.553465835839999
$ws.Cells.Item(4, 19).Value = 0.003440486694232503
$ws.Cells.Item(4, 20).Value = 0.003440486694232503

# Row 5: FAPs -> ECs (Vtn/Tnfrsf11b)
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Vtn"
$ws.Cells.Item(5, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 21.05317333333333
$ws.Cells.Item(5, 8).Value = 63.15952
$ws.Cells.Item(5, 9).Value = 0.3580298485789791
$ws.Cells.Item(5, 10).Value = 0.3580298485789791
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.06861733333333334
$ws.Cells.Item(5, 14).Value = 0.205852
$ws.Cells.Item(5, 15).Value = 0.01654048691795588
$ws.Cells.Item(5, 16).Value = 0.01654048691795588
$ws.Cells.Item(5, 17).Value = 1.444612612337778
$ws.Cells.Item(5, 18).Value = 13.00151351104
$ws.Cells.Item(5, 19).Value = 0.00592198802665833
$ws.Cells.Item(5, 20).Value = 0.00592198802665833

# Row 6: FAPs -> FAPs (Vtn/Tnfrsf11b)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Vtn"
$ws.Cells.Item(6, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 21.05317333333333
$ws.Cells.Item(6, 8).Value = 63.15952
$ws.Cells.Item(6, 9).Value = 0.3580298485789791
$ws.Cells.Item(6, 10).Value = 0.3580298485789791
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.776574666666666
$ws.Cells.Item(6, 14).Value = 11.329724
$ws.Cells.Item(6, 15).Value = 0.9103586635352137
$ws.Cells.Item(6, 16).Value = 0.9103586635352137
$ws.Cells.Item(6, 17).Value = 79.50888106360888
$ws.Cells.Item(6, 18).Value = 715.57992957248
$ws.Cells.Item(6, 19).Value = 0.3259355744580744
$ws.Cells.Item(6, 20).Value = 0.3259355744580744

# Row 7: FAPs -> MuSCs (Vtn/Tnfrsf11b)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Vtn"
$ws.Cells.Item(7, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 21.05317333333333
$ws.Cells.Item(7, 8).Value = 63.15952
$ws.Cells.Item(7, 9).Value = 0.3580298485789791
$ws.Cells.Item(7, 10).Value = 0.3580298485789791
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.303255
$ws.Cells.Item(7, 14).Value = 0.909765
$ws.Cells.Item(7, 15).Value = 0.07310084954683041
$ws.Cells.Item(7, 16).Value = 0.07310084954683042
$ws.Cells.Item(7, 17).Value = 6.3844800792
$ws.Cells.Item(7, 18).Value = 57.46032071280001
$ws.Cells.Item(7, 19).Value = 0.02617228609424643
$ws.Cells.Item(7, 20).Value = 0.02617228609424643

# Row 8: MuSCs -> ECs (Vtn/Tnfrsf11b)
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Vtn"
$ws.Cells.Item(8, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 34.97741266666667
$ws.Cells.Item(8, 8).Value = 104.932238
$ws.Cells.Item(8, 9).Value = 0.5948251867999219
$ws.Cells.Item(8, 10).Value = 0.5948251867999219
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.06861733333333334
$ws.Cells.Item(8, 14).Value = 0.205852
$ws.Cells.Item(8, 15).Value = 0.01654048691795588
$ws.Cells.Item(8, 16).Value = 0.01654048691795588
$ws.Cells.Item(8, 17).Value = 2.400056784086223
$ws.Cells.Item(8, 18).Value = 21.600511056776
$ws.Cells.Item(8, 19).Value = 0.009838698220734774
$ws.Cells.Item(8, 20).Value = 0.009838698220734774

# Row 9: MuSCs -> FAPs (Vtn/Tnfrsf11b)
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Vtn"
$ws.Cells.Item(9, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 34.97741266666667
$ws.Cells.Item(9, 8).Value = 104.932238
$ws.Cells.Item(9, 9).Value = 0.5948251867999219
$ws.Cells.Item(9, 10).Value = 0.5948251867999219
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.776574666666666
$ws.Cells.Item(9, 14).Value = 11.329724
$ws.Cells.Item(9, 15).Value = 0.9103586635352137
$ws.Cells.Item(9, 16).Value = 0.9103586635352137
$ws.Cells.Item(9, 17).Value = 132.0948105824791
$ws.Cells.Item(9, 18).Value = 1188.853295242312
$ws.Cells.Item(9, 19).Value = 0.5415042620922608
$ws.Cells.Item(9, 20).Value = 0.5415042620922608

# Row 10: MuSCs -> MuSCs (Vtn/Tnfrsf11b)
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Vtn"
$ws.Cells.Item(10, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 34.97741266666667
$ws.Cells.Item(10, 8).Value = 104.932238
$ws.Cells.Item(10, 9).Value = 0.5948251867999219
$ws.Cells.Item(10, 10).Value = 0.5948251867999219
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.303255
$ws.Cells.Item(10, 14).Value = 0.909765
$ws.Cells.Item(10, 15).Value = 0.07310084954683041
$ws.Cells.Item(10, 16).Value = 0.07310084954683042
$ws.Cells.Item(10, 17).Value = 10.60707527823
$ws.Cells.Item(10, 18).Value = 95.46367750407002
$ws.Cells.Item(10, 19).Value = 0.04348222648692639
$ws.Cells.Item(10, 20).Value = 0.04348222648692639

# Row 11: Resolving-Mac -> ECs (Vtn/Tnfrsf11b)
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Vtn"
$ws.Cells.Item(11, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.004706
$ws.Cells.Item(11, 8).Value = 0.014118
$ws.Cells.Item(11, 9).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(11, 10).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.06861733333333334
$ws.Cells.Item(11, 14).Value = 0.205852
$ws.Cells.Item(11, 15).Value = 0.01654048691795588
$ws.Cells.Item(11, 16).Value = 0.01654048691795588
$ws.Cells.Item(11, 17).Value = 0.0003229131706666667
$ws.Cells.Item(11, 18).Value = 0.002906218536
$ws.Cells.Item(11, 19).Value = [double]"1.323737529359981E-06"
$ws.Cells.Item(11, 20).Value = [double]"1.323737529359981E-06"

# Row 12: Resolving-Mac -> FAPs (Vtn/Tnfrsf11b)
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Vtn"
$ws.Cells.Item(12, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.004706
$ws.Cells.Item(12, 8).Value = 0.014118
$ws.Cells.Item(12, 9).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(12, 10).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.776574666666666
$ws.Cells.Item(12, 14).Value = 11.329724
$ws.Cells.Item(12, 15).Value = 0.9103586635352137
$ws.Cells.Item(12, 16).Value = 0.9103586635352137
$ws.Cells.Item(12, 17).Value = 0.01777256038133333
$ws.Cells.Item(12, 18).Value = 0.159953043432
$ws.Cells.Item(12, 19).Value = [double]"7.285613380530907E-05"
$ws.Cells.Item(12, 20).Value = [double]"7.285613380530907E-05"

# Row 13: Resolving-Mac -> MuSCs (Vtn/Tnfrsf11b)
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Vtn"
$ws.Cells.Item(13, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.004706
$ws.Cells.Item(13, 8).Value = 0.014118
$ws.Cells.Item(13, 9).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(13, 10).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.303255
$ws.Cells.Item(13, 14).Value = 0.909765
$ws.Cells.Item(13, 15).Value = 0.07310084954683041
$ws.Cells.Item(13, 16).Value = 0.07310084954683042
$ws.Cells.Item(13, 17).Value = 0.00142711803
$ws.Cells.Item(13, 18).Value = 0.01284406227
$ws.Cells.Item(13, 19).Value = [double]"5.850271425092703E-06"
$ws.Cells.Item(13, 20).Value = [double]"5.850271425092704E-06"
